$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Gender" values in column W for rows 4-10 (new shared strings Male/Female)
$ws.Range("W4").Value2 = "Male"
$ws.Range("W5").Value2 = "Male"
$ws.Range("W6").Value2 = "Female"
$ws.Range("W7").Value2 = "Male"
$ws.Range("W8").Value2 = "Male"
$ws.Range("W9").Value2 = "Female"
$ws.Range("W10").Value2 = "Female"

# Update the active selection on the frozen bottom-right pane to follow the new column
$ws.Activate()
$ws.Range("W11").Select()
